$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 0.1931818181818182
$ws.Cells.Item(2, 3).Value = 0.5378787878787878
$ws.Cells.Item(2, 10).Value = 0.01893939393939394
$ws.Cells.Item(2, 16).Value = 0.1363636363636364
$ws.Cells.Item(2, 19).Value = 0.1136363636363636

# Row 3
$ws.Cells.Item(3, 2).Value = 0.006622516556291391
$ws.Cells.Item(3, 3).Value = 0.05960264900662252
$ws.Cells.Item(3, 10).Value = 0.01324503311258278
$ws.Cells.Item(3, 16).Value = 0.6754966887417219
$ws.Cells.Item(3, 19).Value = 0.2450331125827815

# Row 4
$ws.Cells.Item(4, 10).Value = 0.04651162790697674
$ws.Cells.Item(4, 15).Value = 0.02325581395348837
$ws.Cells.Item(4, 16).Value = 0.6744186046511628
$ws.Cells.Item(4, 19).Value = 0.2558139534883721

# Row 6
$ws.Cells.Item(6, 2).Value = 0.04797047970479705
$ws.Cells.Item(6, 4).Value = 0.007380073800738007
$ws.Cells.Item(6, 6).Value = 0.1070110701107011
$ws.Cells.Item(6, 10).Value = 0.2177121771217712
$ws.Cells.Item(6, 15).Value = 0.01476014760147601
$ws.Cells.Item(6, 17).Value = 0.1143911439114391
$ws.Cells.Item(6, 18).Value = 0.1033210332103321
$ws.Cells.Item(6, 19).Value = 0.3874538745387454

# Row 7
$ws.Cells.Item(7, 2).Value = 0.1073446327683616
$ws.Cells.Item(7, 4).Value = 0.02259887005649718
$ws.Cells.Item(7, 5).Value = 0.005649717514124294
$ws.Cells.Item(7, 6).Value = 0.06214689265536723
$ws.Cells.Item(7, 10).Value = 0.1299435028248588
$ws.Cells.Item(7, 15).Value = 0.01129943502824859
$ws.Cells.Item(7, 17).Value = 0.1581920903954802
$ws.Cells.Item(7, 18).Value = 0.1412429378531073
$ws.Cells.Item(7, 19).Value = 0.3615819209039548

# Row 8
$ws.Cells.Item(8, 2).Value = 0.1075
$ws.Cells.Item(8, 4).Value = 0.02
$ws.Cells.Item(8, 5).Value = 0.0025
$ws.Cells.Item(8, 6).Value = 0.05
$ws.Cells.Item(8, 10).Value = 0.1025
$ws.Cells.Item(8, 15).Value = 0.0125
$ws.Cells.Item(8, 17).Value = 0.1525
$ws.Cells.Item(8, 18).Value = 0.1125
$ws.Cells.Item(8, 19).Value = 0.44

# Row 9
$ws.Cells.Item(9, 2).Value = 0.06477732793522267
$ws.Cells.Item(9, 4).Value = 0.01619433198380567
$ws.Cells.Item(9, 5).Value = 0.004048582995951417
$ws.Cells.Item(9, 6).Value = 0.0728744939271255
$ws.Cells.Item(9, 10).Value = 0.07692307692307693
$ws.Cells.Item(9, 15).Value = 0.02834008097165992
$ws.Cells.Item(9, 17).Value = 0.1700404858299595
$ws.Cells.Item(9, 18).Value = 0.09716599190283401
$ws.Cells.Item(9, 19).Value = 0.4696356275303644

# Row 10
$ws.Cells.Item(10, 2).Value = 0.101095197978096
$ws.Cells.Item(10, 4).Value = 0.02190395956192081
$ws.Cells.Item(10, 5).Value = 0.002527379949452401
$ws.Cells.Item(10, 6).Value = 0.09267059814658804
$ws.Cells.Item(10, 10).Value = 0.09435551811288964
$ws.Cells.Item(10, 15).Value = 0.006739679865206402
$ws.Cells.Item(10, 17).Value = 0.2055602358887953
$ws.Cells.Item(10, 18).Value = 0.1069924178601516
$ws.Cells.Item(10, 19).Value = 0.3681550126368998

# Row 11
$ws.Cells.Item(11, 7).Value = 0.1814671814671815
$ws.Cells.Item(11, 10).Value = 0.08108108108108109
$ws.Cells.Item(11, 11).Value = 0.1969111969111969
$ws.Cells.Item(11, 12).Value = 0.5366795366795367
$ws.Cells.Item(11, 19).Value = 0.003861003861003861

# Row 12
$ws.Cells.Item(12, 7).Value = 0.7152777777777778
$ws.Cells.Item(12, 10).Value = 0.2222222222222222
$ws.Cells.Item(12, 11).Value = 0.01388888888888889
$ws.Cells.Item(12, 12).Value = 0.02777777777777778
$ws.Cells.Item(12, 19).Value = 0.02083333333333333

# Row 13
$ws.Cells.Item(13, 6).Value = 0.02040816326530612
$ws.Cells.Item(13, 7).Value = 0.6530612244897959
$ws.Cells.Item(13, 10).Value = 0.3061224489795918
$ws.Cells.Item(13, 19).Value = 0.02040816326530612

# Row 15
$ws.Cells.Item(15, 6).Value = 0.02
$ws.Cells.Item(15, 8).Value = 0.156
$ws.Cells.Item(15, 9).Value = 0.116
$ws.Cells.Item(15, 10).Value = 0.384
$ws.Cells.Item(15, 11).Value = 0.064
$ws.Cells.Item(15, 13).Value = 0.012
$ws.Cells.Item(15, 15).Value = 0.104
$ws.Cells.Item(15, 19).Value = 0.144

# Row 16
$ws.Cells.Item(16, 6).Value = 0.01875
$ws.Cells.Item(16, 9).Value = 0.08749999999999999
$ws.Cells.Item(16, 10).Value = 0.43125
$ws.Cells.Item(16, 11).Value = 0.125
$ws.Cells.Item(16, 15).Value = 0.08125
$ws.Cells.Item(16, 19).Value = 0.13125

# Row 17
$ws.Cells.Item(17, 6).Value = 0.02463054187192118
$ws.Cells.Item(17, 8).Value = 0.1650246305418719
$ws.Cells.Item(17, 9).Value = 0.1133004926108374
$ws.Cells.Item(17, 10).Value = 0.3940886699507389
$ws.Cells.Item(17, 11).Value = 0.07142857142857142
$ws.Cells.Item(17, 13).Value = 0.01724137931034483
$ws.Cells.Item(17, 15).Value = 0.06896551724137931
$ws.Cells.Item(17, 19).Value = 0.145320197044335

# Row 18
$ws.Cells.Item(18, 6).Value = 0.0321285140562249
$ws.Cells.Item(18, 8).Value = 0.2008032128514056
$ws.Cells.Item(18, 9).Value = 0.09236947791164658
$ws.Cells.Item(18, 10).Value = 0.357429718875502
$ws.Cells.Item(18, 11).Value = 0.07630522088353414
$ws.Cells.Item(18, 13).Value = 0.02409638554216868
$ws.Cells.Item(18, 15).Value = 0.1124497991967871
$ws.Cells.Item(18, 19).Value = 0.1044176706827309

# Row 19
$ws.Cells.Item(19, 6).Value = 0.01537216828478964
$ws.Cells.Item(19, 8).Value = 0.1868932038834951
$ws.Cells.Item(19, 9).Value = 0.1108414239482201
$ws.Cells.Item(19, 10).Value = 0.3689320388349515
$ws.Cells.Item(19, 11).Value = 0.09304207119741101
$ws.Cells.Item(19, 13).Value = 0.02427184466019417
$ws.Cells.Item(19, 14).Value = 0.0008090614886731392
$ws.Cells.Item(19, 15).Value = 0.09385113268608414
$ws.Cells.Item(19, 19).Value = 0.1059870550161812
